# Fruta / hortaliza, semanal
# Rewrites rows 3-13 of the Aji sheet with the updated weekly values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44421
$ws.Range("H3").Value = 'Americana (o)'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 75000
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = 75000
$ws.Range("N3").Value = '$/caja 25 kilos'
$ws.Range("P3").Value = 3000
$ws.Range("Q3").Value = 25

# Row 4
$ws.Range("D4").Value = 44319
$ws.Range("H4").Value = 'Americana (o)'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 30000
$ws.Range("N4").Value = '$/caja 25 kilos'
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = 25

# Row 5
$ws.Range("D5").Value = 44449
$ws.Range("H5").Value = 'Americana (o)'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 80000
$ws.Range("L5").Value = 80000
$ws.Range("M5").Value = 80000
$ws.Range("N5").Value = '$/caja 25 kilos'
$ws.Range("P5").Value = 3200
$ws.Range("Q5").Value = 25

# Row 6
$ws.Range("D6").Value = 44449
$ws.Range("H6").Value = 'Americana (o)'
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 75000
$ws.Range("L6").Value = 75000
$ws.Range("M6").Value = 75000
$ws.Range("N6").Value = '$/caja 15 kilos'
$ws.Range("P6").Value = 5000
$ws.Range("Q6").Value = 15

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("H7").Value = 'Americana (o)'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 78000
$ws.Range("L7").Value = 78000
$ws.Range("M7").Value = 78000
$ws.Range("N7").Value = '$/caja 25 kilos'
$ws.Range("P7").Value = 3120
$ws.Range("Q7").Value = 25

# Row 8
$ws.Range("D8").Value = 44446
$ws.Range("H8").Value = 'Inferno'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 80000
$ws.Range("L8").Value = 80000
$ws.Range("M8").Value = 80000
$ws.Range("N8").Value = '$/caja 15 kilos'
$ws.Range("P8").Value = 5333
$ws.Range("Q8").Value = 15

# Row 9
$ws.Range("D9").Value = 44340
$ws.Range("H9").Value = 'Americana (o)'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 35000
$ws.Range("L9").Value = 35000
$ws.Range("M9").Value = 35000
$ws.Range("N9").Value = '$/caja 25 kilos'
$ws.Range("P9").Value = 1400
$ws.Range("Q9").Value = 25

# Row 10
$ws.Range("D10").Value = 44193
$ws.Range("H10").Value = 'Americana (o)'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 46000
$ws.Range("L10").Value = 46000
$ws.Range("M10").Value = 46000
$ws.Range("N10").Value = '$/caja 15 kilos'
$ws.Range("P10").Value = 3067
$ws.Range("Q10").Value = 15

# Row 11
$ws.Range("D11").Value = 44343
$ws.Range("H11").Value = 'Americana (o)'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 36000
$ws.Range("L11").Value = 36000
$ws.Range("M11").Value = 36000
$ws.Range("N11").Value = '$/caja 25 kilos'
$ws.Range("P11").Value = 1440
$ws.Range("Q11").Value = 25

# Row 12
$ws.Range("D12").Value = 44221
$ws.Range("H12").Value = 'Americana (o)'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 22
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 24545
$ws.Range("N12").Value = '$/caja 25 kilos'
$ws.Range("P12").Value = 982
$ws.Range("Q12").Value = 25

# Row 13
$ws.Range("D13").Value = 44425
$ws.Range("H13").Value = 'Americana (o)'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 15
$ws.Range("K13").Value = 75000
$ws.Range("L13").Value = 75000
$ws.Range("M13").Value = 75000
$ws.Range("N13").Value = '$/caja 25 kilos'
$ws.Range("P13").Value = 3000
$ws.Range("Q13").Value = 25

